# This script applies the weekly Fruta/Hortaliza price-log update for
# "Femacal de La Calera - Platano": the most recent 10 weekly rows (423-432)
# shift down by one reporting period (their old contents move to rows
# 425-434), and a brand-new latest week (date 2021-09-09 / serial 44448)
# is written into rows 423-424.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Roll the existing weekly blocks (rows 423-432) down by one period ---
# Row 423
$ws.Range("D423").Value = 44448
$ws.Range("L423").Value = "Pintón"
$ws.Range("M423").Value = 360
$ws.Range("N423").Value = 18000
$ws.Range("O423").Value = 19000
$ws.Range("P423").Value = 18444
$ws.Range("S423").Value = 922

# Row 424
$ws.Range("D424").Value = 44448
$ws.Range("L424").Value = "Primera Pintón"
$ws.Range("M424").Value = 160
$ws.Range("N424").Value = 20000
$ws.Range("O424").Value = 20000
$ws.Range("P424").Value = 20000
$ws.Range("S424").Value = 1000

# Row 425
$ws.Range("D425").Value = 44167
$ws.Range("L425").Value = "Pintón"
$ws.Range("M425").Value = 360
$ws.Range("N425").Value = 17000
$ws.Range("O425").Value = 17000
$ws.Range("P425").Value = 17000
$ws.Range("S425").Value = 850

# Row 426
$ws.Range("D426").Value = 44167
$ws.Range("L426").Value = "Primera Pintón"
$ws.Range("M426").Value = 240
$ws.Range("N426").Value = 19000
$ws.Range("O426").Value = 19000
$ws.Range("P426").Value = 19000
$ws.Range("S426").Value = 950

# Row 427
$ws.Range("D427").Value = 44238
$ws.Range("L427").Value = "Maduro"
$ws.Range("M427").Value = 250
$ws.Range("N427").Value = 11000
$ws.Range("O427").Value = 11000
$ws.Range("P427").Value = 11000
$ws.Range("S427").Value = 550

# Row 428
$ws.Range("D428").Value = 44238
$ws.Range("L428").Value = "Pintón"
$ws.Range("M428").Value = 280
$ws.Range("N428").Value = 12000
$ws.Range("O428").Value = 12000
$ws.Range("P428").Value = 12000
$ws.Range("S428").Value = 600

# Row 429
$ws.Range("D429").Value = 44238
$ws.Range("L429").Value = "Verde"
$ws.Range("M429").Value = 180
$ws.Range("N429").Value = 11500
$ws.Range("O429").Value = 11500
$ws.Range("P429").Value = 11500
$ws.Range("S429").Value = 575

# Row 430
$ws.Range("D430").Value = 44399
$ws.Range("L430").Value = "Maduro"
$ws.Range("M430").Value = 160
$ws.Range("N430").Value = 17000
$ws.Range("O430").Value = 17000
$ws.Range("P430").Value = 17000
$ws.Range("S430").Value = 850

# Row 431
$ws.Range("D431").Value = 44399
$ws.Range("L431").Value = "Pintón"
$ws.Range("M431").Value = 320
$ws.Range("N431").Value = 18000
$ws.Range("O431").Value = 18000
$ws.Range("P431").Value = 18000
$ws.Range("S431").Value = 900

# Row 432
$ws.Range("D432").Value = 44399
$ws.Range("L432").Value = "Primera Pintón"
$ws.Range("M432").Value = 280
$ws.Range("N432").Value = 20000
$ws.Range("O432").Value = 20000
$ws.Range("P432").Value = 20000
$ws.Range("S432").Value = 1000

# --- Append the two newly-reported rows (433-434) ---
# Row 433
$ws.Range("A433").Value = 3
$ws.Range("B433").Value = "Femacal de La Calera"
$ws.Range("C433").Value = "Coquimbo"
$ws.Range("D433").Value = 44400
$ws.Range("E433").Value = 5
$ws.Range("F433").Value = "Fruta"
$ws.Range("G433").Value = 100108
$ws.Range("H433").Value = "Tropicales y subtropicales"
$ws.Range("I433").Value = 100108006
$ws.Range("J433").Value = "Plátano"
$ws.Range("K433").Value = "Sin especificar"
$ws.Range("L433").Value = "Pintón"
$ws.Range("M433").Value = 320
$ws.Range("N433").Value = 18000
$ws.Range("O433").Value = 18000
$ws.Range("P433").Value = 18000
$ws.Range("Q433").Value = "$/caja 20 kilos"
$ws.Range("R433").Value = "Ecuador"
$ws.Range("S433").Value = 900
$ws.Range("T433").Value = 20
$ws.Range("D433").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 434
$ws.Range("A434").Value = 3
$ws.Range("B434").Value = "Femacal de La Calera"
$ws.Range("C434").Value = "Coquimbo"
$ws.Range("D434").Value = 44400
$ws.Range("E434").Value = 5
$ws.Range("F434").Value = "Fruta"
$ws.Range("G434").Value = 100108
$ws.Range("H434").Value = "Tropicales y subtropicales"
$ws.Range("I434").Value = 100108006
$ws.Range("J434").Value = "Plátano"
$ws.Range("K434").Value = "Sin especificar"
$ws.Range("L434").Value = "Primera Pintón"
$ws.Range("M434").Value = 280
$ws.Range("N434").Value = 20000
$ws.Range("O434").Value = 20000
$ws.Range("P434").Value = 20000
$ws.Range("Q434").Value = "$/caja 20 kilos"
$ws.Range("R434").Value = "Ecuador"
$ws.Range("S434").Value = 1000
$ws.Range("T434").Value = 20
$ws.Range("D434").NumberFormat = "YYYY-MM-DD HH:MM:SS"

